# Update cryptos.xlsx data (Thu Apr  6 23:45:28 UTC 2023 refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "28.088.61" }
    @{ Cell = "D3"; Value = "1.875.14" }
    @{ Cell = "D4"; Value = "1.002" }
    @{ Cell = "E4"; Value = "  +0.18%  " }
    @{ Cell = "D5"; Value = "313.19" }
    @{ Cell = "E5"; Value = "  -0.39%  " }
    @{ Cell = "D7"; Value = "0.5048" }
    @{ Cell = "E7"; Value = "  -0.27%  " }
    @{ Cell = "D8"; Value = "0.3842" }
    @{ Cell = "E8"; Value = "  -2.22%  " }
    @{ Cell = "D9"; Value = "0.08600" }
    @{ Cell = "E9"; Value = "  -8.16%  " }
    @{ Cell = "D10"; Value = "1.116" }
    @{ Cell = "E10"; Value = "  -2.44%  " }
    @{ Cell = "E11"; Value = "  -1.76%  " }
    @{ Cell = "E12"; Value = "  -1.87%  " }
    @{ Cell = "B13"; Value = "Solana" }
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol" }
    @{ Cell = "D13"; Value = "20.66" }
    @{ Cell = "E13"; Value = "  -1.78%  " }
    @{ Cell = "B14"; Value = "WrappedEther" }
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" }
    @{ Cell = "D14"; Value = "1.873.93" }
    @{ Cell = "E14"; Value = "  -2.09%  " }
    @{ Cell = "D15"; Value = "7.212" }
    @{ Cell = "E15"; Value = "  -1.64%  " }
    @{ Cell = "D16"; Value = "1.002" }
    @{ Cell = "E16"; Value = "  +0.17%  " }
    @{ Cell = "D17"; Value = "0.00001099" }
    @{ Cell = "E17"; Value = "  -2.53%  " }
    @{ Cell = "D18"; Value = "91.10" }
    @{ Cell = "E18"; Value = "  -1.69%  " }
    @{ Cell = "D19"; Value = "0.06627" }
    @{ Cell = "E19"; Value = "  +0.06%  " }
    @{ Cell = "E20"; Value = "  +0.28%  " }
    @{ Cell = "D22"; Value = "6.099" }
    @{ Cell = "E22"; Value = "  -2.23%  " }
    @{ Cell = "D23"; Value = "28.124.57" }
    @{ Cell = "E23"; Value = "  -0.48%  " }
    @{ Cell = "E24"; Value = "  -1.66%  " }
    @{ Cell = "D25"; Value = "2.268" }
    @{ Cell = "E25"; Value = "  -2.41%  " }
    @{ Cell = "B26"; Value = "LEO" }
    @{ Cell = "C26"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" }
    @{ Cell = "D26"; Value = "3.405" }
    @{ Cell = "E26"; Value = "  +0.42%  " }
    @{ Cell = "B27"; Value = "LidoDAOToken" }
    @{ Cell = "C27"; Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo" }
    @{ Cell = "D27"; Value = "2.592" }
    @{ Cell = "E27"; Value = "  -0.36%  " }
    @{ Cell = "B28"; Value = "WrappedliquidstakedEther2.0" }
    @{ Cell = "C28"; Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth" }
    @{ Cell = "D28"; Value = "2.091.38" }
    @{ Cell = "E28"; Value = "  -2.06%  " }
    @{ Cell = "B29"; Value = "EthereumClassic" }
    @{ Cell = "C29"; Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc" }
    @{ Cell = "D29"; Value = "20.72" }
    @{ Cell = "E29"; Value = "  -2.20%  " }
    @{ Cell = "B30"; Value = "Monero" }
    @{ Cell = "C30"; Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" }
    @{ Cell = "D30"; Value = "156.85" }
    @{ Cell = "E30"; Value = "  -0.89%  " }
    @{ Cell = "B31"; Value = "BitcoinCash" }
    @{ Cell = "C31"; Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch" }
    @{ Cell = "D31"; Value = "126.45" }
    @{ Cell = "E31"; Value = "  -0.71%  " }
    @{ Cell = "B32"; Value = "Stellar" }
    @{ Cell = "C32"; Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm" }
    @{ Cell = "D32"; Value = "0.1058" }
    @{ Cell = "E32"; Value = "  -1.35%  " }
    @{ Cell = "B33"; Value = "ImmutableX" }
    @{ Cell = "C33"; Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" }
    @{ Cell = "D33"; Value = "1.062" }
    @{ Cell = "E33"; Value = "  -4.16%  " }
    @{ Cell = "B34"; Value = "Filecoin" }
    @{ Cell = "C34"; Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil" }
    @{ Cell = "D34"; Value = "5.617" }
    @{ Cell = "E34"; Value = "  -1.00%  " }
    @{ Cell = "B35"; Value = "HuobiToken" }
    @{ Cell = "C35"; Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht" }
    @{ Cell = "D35"; Value = "3.588" }
    @{ Cell = "E35"; Value = "  -0.53%  " }
    @{ Cell = "B36"; Value = "FraxShare" }
    @{ Cell = "C36"; Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs" }
    @{ Cell = "D36"; Value = "9.638" }
    @{ Cell = "E36"; Value = "  -0.66%  " }
    @{ Cell = "B37"; Value = "VeChain" }
    @{ Cell = "C37"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet" }
    @{ Cell = "D37"; Value = "0.02454" }
    @{ Cell = "E37"; Value = "  +0.42%  " }
    @{ Cell = "B38"; Value = "Hedera" }
    @{ Cell = "C38"; Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar" }
    @{ Cell = "D38"; Value = "0.06582" }
    @{ Cell = "E38"; Value = "  -1.99%  " }
    @{ Cell = "B39"; Value = "Algorand" }
    @{ Cell = "C39"; Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo" }
    @{ Cell = "D39"; Value = "0.2183" }
    @{ Cell = "E39"; Value = "  -1.65%  " }
    @{ Cell = "B40"; Value = "ARBITRUM" }
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb" }
    @{ Cell = "D40"; Value = "1.213" }
    @{ Cell = "E40"; Value = "  -2.87%  " }
    @{ Cell = "B41"; Value = "TrustWalletToken" }
    @{ Cell = "C41"; Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt" }
    @{ Cell = "D41"; Value = "1.243" }
    @{ Cell = "E41"; Value = "  -2.95%  " }
    @{ Cell = "B42"; Value = "TheSandbox" }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand" }
    @{ Cell = "D42"; Value = "0.6386" }
    @{ Cell = "E42"; Value = "  -2.51%  " }
    @{ Cell = "B43"; Value = "Aptos" }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt" }
    @{ Cell = "D43"; Value = "11.45" }
    @{ Cell = "E43"; Value = "  -1.01%  " }
    @{ Cell = "B44"; Value = "InternetComputer(DFINITY)" }
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp" }
    @{ Cell = "D44"; Value = "4.901" }
    @{ Cell = "E44"; Value = "  -2.60%  " }
    @{ Cell = "B45"; Value = "EnergySwap" }
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" }
    @{ Cell = "D45"; Value = "13.22" }
    @{ Cell = "E45"; Value = "  -1.74%  " }
    @{ Cell = "B46"; Value = "Decentraland" }
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana" }
    @{ Cell = "D46"; Value = "0.6010" }
    @{ Cell = "E46"; Value = "  -2.12%  " }
    @{ Cell = "B47"; Value = "WEMIXTOKEN" }
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix" }
    @{ Cell = "D47"; Value = "1.283" }
    @{ Cell = "E47"; Value = "  -1.73%  " }
    @{ Cell = "B48"; Value = "PancakeSwap" }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake" }
    @{ Cell = "D48"; Value = "3.665" }
    @{ Cell = "E48"; Value = "  -1.58%  " }
    @{ Cell = "B49"; Value = "NEARProtocol" }
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near" }
    @{ Cell = "D49"; Value = "1.992" }
    @{ Cell = "E49"; Value = "  -1.97%  " }
    @{ Cell = "B50"; Value = "EOS" }
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos" }
    @{ Cell = "D50"; Value = "1.223" }
    @{ Cell = "E50"; Value = "  +2.72%  " }
    @{ Cell = "B51"; Value = "Quant" }
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt" }
    @{ Cell = "D51"; Value = "121.45" }
    @{ Cell = "E51"; Value = "  -0.86%  " }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}

Write-Host "Applied cryptos.xlsx update: $($updates.Count) cells updated"
